$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.177.91'
$ws.Range("E2").Value = '  +2.13%  '

$ws.Range("D3").Value = '3.176.22'
$ws.Range("E3").Value = '  +4.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.19'
$ws.Range("E5").Value = '  +4.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.84'
$ws.Range("E6").Value = '  +7.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.175.48'
$ws.Range("E8").Value = '  +4.13%  '

$ws.Range("E9").Value = '  +2.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  +6.52%  '

$ws.Range("E11").Value = '  -0.09%  '

$ws.Range("E12").Value = '  +3.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  +19.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.71'
$ws.Range("E14").Value = '  +7.02%  '

$ws.Range("D15").Value = '3.694.71'
$ws.Range("E15").Value = '  +4.29%  '

$ws.Range("D16").Value = '65.218.46'
$ws.Range("E16").Value = '  +2.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  +6.45%  '

$ws.Range("D18").Value = '3.170.91'
$ws.Range("E18").Value = '  +4.15%  '

$ws.Range("E19").Value = '  +1.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '514.27'
$ws.Range("E20").Value = '  +6.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.91'
$ws.Range("E21").Value = '  +5.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.729'
$ws.Range("E22").Value = '  +6.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.49'
$ws.Range("E23").Value = '  +6.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.86'
$ws.Range("E24").Value = '  +4.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.22'
$ws.Range("E25").Value = '  +3.30%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.13'
$ws.Range("E27").Value = '  +12.57%  '

$ws.Range("E28").Value = '  +5.27%  '

$ws.Range("E29").Value = '  +8.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.93'
$ws.Range("E30").Value = '  +6.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.83'
$ws.Range("E31").Value = '  +16.13%  '

$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("E33").Value = '  +4.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.34'
$ws.Range("E34").Value = '  +11.98%  '

$ws.Range("E35").Value = '  +6.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.90'
$ws.Range("E36").Value = '  +1.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0905'
$ws.Range("E37").Value = '  +11.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '478.44'
$ws.Range("E38").Value = '  +8.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.12'
$ws.Range("E39").Value = '  +12.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0423'

$ws.Range("E41").Value = '  +4.83%  '

$ws.Range("D42").Value = '3.085.89'
$ws.Range("E42").Value = '  +2.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.119'
$ws.Range("E43").Value = '  +2.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.47'
$ws.Range("E44").Value = '  +10.78%  '

$ws.Range("E45").Value = '  +5.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.43'
$ws.Range("E46").Value = '  +6.11%  '

$ws.Range("D47").Value = '0.0₃0609'
$ws.Range("E47").Value = '  +19.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("E49").Value = '  +1.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.27'
$ws.Range("E50").Value = '  +8.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.57'
$ws.Range("E51").Value = '  +2.11%  '
